# Update "Hjemme passive" data - meanEMG legmaxROM values (row 1)
# and recompute CON / STR rows (2 and 3) for columns B:E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - updated leg max ROM header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) - updated meanEMG values; C2 and E2 no longer have data
$ws.Range("B2").Value = 24.680170421538492
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 15.621528783809993
$ws.Range("E2").ClearContents()

# Row 3 (STR) - updated meanEMG values
$ws.Range("B3").Value = 21.358752986927641
$ws.Range("C3").Value = 22.112663767150138
$ws.Range("D3").Value = 15.794170038206403
$ws.Range("E3").Value = 33.035550857034877

# Selection now only spans the updated data range
$ws.Range("B1:E3").Select()
